$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 533.3333
$ws.Range("I2").Value = 533.3333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 533.3333
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -420.3333
$ws.Range("N2").ClearContents()

$ws.Range("H28").Value = 276.33334
$ws.Range("I28").Value = 185.875
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 185.875
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 299.125
$ws.Range("N28").Value = -1970

$ws.Range("H123").Value = 48000
$ws.Range("J123").Value = 48000
$ws.Range("L123").Value = 48000
$ws.Range("N123").Value = -57800

$ws.Range("H137").Value = 1836.35
$ws.Range("I137").Value = 1513.3529
$ws.Range("K137").Value = 4540.0587
$ws.Range("M137").Value = -1990.0587

$ws.Range("H138").Value = 1963.0968
$ws.Range("J138").Value = 2242.7856
$ws.Range("L138").Value = 6728.3568
$ws.Range("N138").Value = -17008.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4410.5425
$ws.Range("I32").Value = 3110.2693
$ws.Range("K32").Value = 3110.2693
$ws.Range("M32").Value = -2823.2693

$ws.Range("H74").Value = 1905.3572
$ws.Range("I74").Value = 854.5714
$ws.Range("K74").Value = 854.5714
$ws.Range("M74").Value = 19.42859999999996

$ws.Range("H77").Value = 1905.3572
$ws.Range("I77").Value = 854.5714
$ws.Range("K77").Value = 4272.857
$ws.Range("M77").Value = 95.14300000000003

$ws.Range("H88").Value = 4233.3335

$ws.Range("H91").Value = 4233.3335

$ws.Range("H122").Value = 1264
$ws.Range("I122").Value = 948.8570999999999
$ws.Range("K122").Value = 2846.5713
$ws.Range("M122").Value = -396.5712999999996

$ws.Range("H132").Value = 1488.3889
$ws.Range("I132").Value = 1502.4117
$ws.Range("K132").Value = 4507.2351
$ws.Range("M132").Value = -1977.2351

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1094385.9
$ws.Range("I86").Value = 1115194
$ws.Range("K86").Value = 1115194
$ws.Range("M86").Value = -1114071

$ws.Range("H89").Value = 1094385.9
$ws.Range("I89").Value = 1115194
$ws.Range("K89").Value = 5575970
$ws.Range("M89").Value = -5570354

$ws.Range("H94").Value = 1535.7
$ws.Range("I94").Value = 976.1667
$ws.Range("J94").Value = 2375
$ws.Range("K94").Value = 976.1667
$ws.Range("L94").Value = 2375
$ws.Range("M94").Value = -525.1667
$ws.Range("N94").Value = -3277

$ws.Range("H134").Value = 6251.364
$ws.Range("I134").Value = 7118.722
$ws.Range("J134").Value = 2348.25
$ws.Range("K134").Value = 21356.166
$ws.Range("L134").Value = 7044.75
$ws.Range("M134").Value = -18821.166
$ws.Range("N134").Value = -12114.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 755.35297
$ws.Range("I107").Value = 534.2143
$ws.Range("J107").Value = 1787.3334
$ws.Range("K107").Value = 534.2143
$ws.Range("L107").Value = 1787.3334
$ws.Range("M107").Value = 1385.7857
$ws.Range("N107").Value = -5627.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 166.125
$ws.Range("I23").Value = 206
$ws.Range("K23").Value = 618
$ws.Range("M23").Value = -383

$ws.Range("H34").Value = 1090.1111
$ws.Range("I34").Value = 96.666664
$ws.Range("J34").Value = 1586.8334
$ws.Range("K34").Value = 289.999992
$ws.Range("L34").Value = 4760.5002
$ws.Range("M34").Value = -205.999992
$ws.Range("N34").Value = -4928.5002

$ws.Range("H68").Value = 1100
$ws.Range("J68").Value = 1100
$ws.Range("L68").Value = 3300
$ws.Range("N68").Value = -4922

$ws.Range("H69").Value = 2825.1667
$ws.Range("I69").Value = 2399.8
$ws.Range("J69").Value = 2988.7693
$ws.Range("K69").Value = 7199.400000000001
$ws.Range("L69").Value = 8966.3079
$ws.Range("M69").Value = -6388.400000000001
$ws.Range("N69").Value = -10588.3079

$ws.Range("H71").Value = 1100
$ws.Range("J71").Value = 1100
$ws.Range("L71").Value = 9900
$ws.Range("N71").Value = -18012

$ws.Range("H72").Value = 2825.1667
$ws.Range("I72").Value = 2399.8
$ws.Range("J72").Value = 2988.7693
$ws.Range("K72").Value = 21598.2
$ws.Range("L72").Value = 26898.9237
$ws.Range("M72").Value = -17542.2
$ws.Range("N72").Value = -35010.9237

$ws.Range("H80").Value = 2355.7144
$ws.Range("J80").Value = 2498.3333
$ws.Range("L80").Value = 7494.999899999999
$ws.Range("N80").Value = -9366.999899999999

$ws.Range("H83").Value = 2355.7144
$ws.Range("J83").Value = 2498.3333
$ws.Range("L83").Value = 22484.9997
$ws.Range("N83").Value = -31844.9997

$ws.Range("H86").Value = 344.5
$ws.Range("J86").Value = 344.5
$ws.Range("L86").Value = 1033.5
$ws.Range("N86").Value = -3405.5

$ws.Range("H89").Value = 344.5
$ws.Range("J89").Value = 344.5
$ws.Range("L89").Value = 3100.5
$ws.Range("N89").Value = -14956.5

$ws.Range("H116").Value = 2447.818
$ws.Range("I116").Value = 1033
$ws.Range("K116").Value = 3099
$ws.Range("M116").Value = 343

$ws.Range("H131").Value = 14515.68
$ws.Range("J131").Value = 15054.922
$ws.Range("L131").Value = 45164.766
$ws.Range("N131").Value = -55244.766

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1105.909
$ws.Range("I113").Value = 963.25
$ws.Range("J113").Value = 1187.4286
$ws.Range("K113").Value = 963.25
$ws.Range("L113").Value = 1187.4286
$ws.Range("M113").Value = 1206.75
$ws.Range("N113").Value = -5527.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2284
$ws.Range("J46").Value = 3155.8
$ws.Range("L46").Value = 3155.8
$ws.Range("N46").Value = -3531.8

$ws.Range("H55").Value = 422
$ws.Range("I55").Value = 336.375
$ws.Range("J55").Value = 696
$ws.Range("K55").Value = 336.375
$ws.Range("L55").Value = 696
$ws.Range("M55").Value = -163.375
$ws.Range("N55").Value = -1042

$ws.Range("H93").Value = 15152441
$ws.Range("I93").Value = 780.82355
$ws.Range("J93").Value = 66668090
$ws.Range("K93").Value = 780.82355
$ws.Range("L93").Value = 66668090
$ws.Range("M93").Value = 467.17645
$ws.Range("N93").Value = -66670586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1106.8
$ws.Range("I100").Value = 933.5
$ws.Range("K100").Value = 1867
$ws.Range("M100").Value = -1326

$ws.Range("H109").Value = 78998.664
$ws.Range("J109").Value = 78998.664
$ws.Range("L109").Value = 78998.664
$ws.Range("N109").Value = -81772.664
